$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-28"

# Update the header label for the current-year column
$ws.Range("I1").Value = "2022 (through 04-28)"

# Update April total (I5) and grand Total (I14) with the new data
$ws.Range("I5").Value = 113
$ws.Range("I14").Value = 548
